# SectorGroup.xlsx (fr) — fix column order so that the "category-name" and
# "group-code" values that were swapped end up in the correct columns.
#
# The workbook has header row 1:
#   A=code  B=name  C=status  D=category-code  E=category-name
#   F=group-code  G=group-name
#
# For every row (including the header, whose labels were also swapped),
# the content that lives in column E actually belongs in column F and
# vice-versa. This script swaps the values of columns E and F for the
# whole used range (rows 1-235) while preserving the original "shared
# string" (text) cell type, using Copy/Paste through a scratch column so
# that Excel does not reinterpret numeric-looking codes (e.g. "110",
# "998") as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$eRange = $ws.Range("E1:E$lastRow")
$fRange = $ws.Range("F1:F$lastRow")
$hRange = $ws.Range("H1:H$lastRow")

# 3-way swap using a scratch column so that copying preserves the
# original cell type (text) instead of Excel inferring a number type
# for values such as "110" or "998".
$eRange.Copy($hRange)
$fRange.Copy($eRange)
$hRange.Copy($fRange)

$hRange.ClearContents()
$excel.CutCopyMode = $false
